$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (price + 1h volume change) pulled on Mon Jan 23 14:56:34 UTC 2023
$updates = @{
    "D2" = "305.27"
    "E2" = "-0.07%"
    "D3" = "35.71"
    "D4" = "5.082"
    "E4" = "1.37%"
    "D5" = "0.07959"
    "E5" = "0.72%"
    "D6" = "2.123"
    "E6" = "-4.19%"
    "D7" = "7.913"
    "E7" = "-1.38%"
    "D8" = "0.9239"
    "E8" = "-0.04%"
    "D9" = "0.09636"
    "E9" = "-0.75%"
    "D10" = "0.1845"
    "E10" = "-2.15%"
    "D11" = "0.08697"
    "E11" = "0.88%"
    "D12" = "0.03547"
    "E12" = "-4.03%"
    "D13" = "0.09909"
    "E13" = "-0.70%"
    "D14" = "0.001431"
    "E14" = "-2.60%"
    "D15" = "0.005659"
    "E15" = "0.17%"
    "D16" = "3.470"
    "E16" = "0.08%"
    "D17" = "4.111"
    "E17" = "2.17%"
    "D18" = "2.753"
    "E18" = "22.40%"
    "E19" = "-1.28%"
    "E20" = "1.55%"
    "D21" = "5.152"
    "E21" = "8.33%"
    "E22" = "0.41%"
    "D23" = "0.04517"
    "E23" = "-1.01%"
    "D24" = "0.001228"
    "E24" = "-0.44%"
    "D25" = "0.004887"
    "E25" = "9.24%"
    "D26" = "0.0001297"
    "E26" = "-7.13%"
    "D27" = "0.0004738"
    "E27" = "-0.28%"
    "D39" = "0.01859"
    "E39" = "0.78%"
    "D40" = "0.04750"
    "E40" = "-1.11%"
    "D41" = "0.007802"
    "E41" = "-3.97%"
    "D42" = "0.1402"
    "E42" = "0.19%"
    "D43" = "0.007717"
    "E43" = "2.21%"
    "D44" = "0.002184"
    "E44" = "-0.94%"
    "D45" = "0.01117"
    "E45" = "11.05%"
    "D46" = "0.00006269"
    "E46" = "0.16%"
    "D47" = "0.00000000748"
    "E47" = "-0.28%"
    "E48" = "0.17%"
    "D49" = "50.65"
    "E49" = "77.00%"
    "D50" = "0.001895"
    "E50" = "10.08%"
    "D51" = "0.00002095"
    "E51" = "-0.28%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to store the value as literal text
    # (matches the source sheet's inline-string cells instead of auto-converting
    # numeric-looking strings like '305.27' or '-0.07%' into Number/Percent types).
    $range.Value = "'" + $updates[$cellRef]
    # Drop the quote-prefix formatting footprint so the cell style is untouched.
    $range.ClearFormats()
}
